# Update "Last Updated" timestamp on the Metadata sheet
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 09:31 AM"

# Update the "% Current Qtr" (F) column on the Industry Analysis sheet
$industry = $wb.Worksheets.Item("Industry Analysis")
$industryF = @{
  2 = 18.476
  3 = -7.7404
  4 = 30.7972
  5 = -50.2266
  6 = 61.9649
  7 = -9.1713
  8 = -3.556
  9 = 38.3509
  10 = -6.2497
  11 = 52.6723
  12 = -6.932
  13 = 17.5662
  14 = -35.5106
  15 = 0.6286
  16 = -3.1514
  17 = -20.6354
  18 = -0.0175
  19 = -26.9255
  20 = 44.703
  21 = 10.0506
  22 = 84.6016
  23 = -54.4868
  24 = -12.8122
  25 = -9.182700000000001
  26 = 5.9529
  27 = -33.2998
  28 = -20.4441
  29 = -17.1514
  30 = 24.527
  31 = 57.6193
  32 = -1.527
  33 = -5.2378
  34 = 27.4054
  35 = 6.7961
  36 = -5.6683
  37 = 1.4178
  38 = -22.4272
  39 = 12.3741
  40 = -5.138
  41 = -0.1825
  42 = 23.2483
  43 = 14.456
  44 = -11.1739
  45 = 27.112
  46 = -5.6252
  47 = -36.5148
  48 = -27.8397
  49 = -25.4424
  50 = -49.1173
  51 = -51.065
  52 = -35.4517
  53 = -11.9879
  54 = -3.0992
  55 = -15.3441
  56 = -25.937
  57 = -29.1486
  58 = -6.4093
  59 = -23.3046
  60 = -11.2657
  61 = -9.777699999999999
  62 = -16.0561
  63 = -9.932499999999999
  64 = 51.8767
  65 = -43.5191
  66 = 13.7315
  67 = 12.6111
  68 = 31.7532
  69 = -19.9577
  70 = -12.9642
  71 = 13.2432
  72 = 2.8232
  73 = -9.179
  74 = -14.2931
  75 = 28.3699
  76 = 45.5868
}
foreach ($r in $industryF.Keys) {
  $industry.Cells.Item([int]$r, 6).Value = $industryF[$r]
}

# Insert a new top row on the Stock List sheet for the newly listed CAPTRU-RE1,
# pushing all existing rows down by one and dropping the last (overflow) row.
$stockList = $wb.Worksheets.Item("Stock List")
$stockList.Rows.Item(2).Insert()
$stockList.Range("A2").Value = "📋"
$stockList.Range("B2").Value = "CAPTRU-RE1"
$stockList.Range("C2").Value = "CAPTRU-RE1"
$stockList.Range("D2").Value = 5.67
$stockList.Range("E2").Value = -11.9565
$stockList.Range("F2").Value = "N/A"
$stockList.Range("G2").Value = "N/A"
$stockList.Range("H2").Value = 0
$stockList.Rows.Item(77).Delete()
